# 1st class week #3
#
# - Bump the `digits = 15` argument to `digits = 16` in Exercicio 2's
#   first code block.
# - Update the printed result of `log(4)` to the longer (16-digit)
#   representation.
# - Add a new pair of source/output blocks showing the same value
#   obtained via sprintf("%.15f", log(4)).

$d = $word.ActiveDocument

# --- 1. "digits = 15" -> "digits = 16" ------------------------------------
$d.Content.Find.Execute("15", $true, $false, $false, $false, $false, $true, `
                         1, $false, "16", 2) | Out-Null

# --- 2. "## [1] 1.38629436111989" -> "## [1] 1.386294361119891" -----------
$d.Content.Find.Execute("## [1] 1.38629436111989", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "## [1] 1.386294361119891", 2) | Out-Null

# --- 3. Insert the two new "SourceCode" paragraphs right after the -------
#        "## [1] 1.386294361119891" paragraph (and before "Exercicio 3").
$hits = $d.Paragraphs
for ($i = 1; $i -le $hits.Count; $i++) {
    $p = $hits.Item($i)
    if ($p.Range.Text -like "*1.386294361119891*") {
        $resultPar = $p
        break
    }
}

# Two fresh empty paragraphs, cloned ("Source Code" style) from $resultPar.
$resultPar.Range.InsertParagraphAfter()
$codePar = $resultPar.Next()
$codePar.Range.InsertParagraphAfter()
$outPar = $codePar.Next()

# --- paragraph 1: sprintf("%.15f", log(4)) --------------------------------
$t = $codePar.Range
$t.Collapse(1)

$t.InsertAfter("sprintf")
$t.Style = "KeywordTok"
$t.Collapse(0)

$t.InsertAfter("(")
$t.Style = "NormalTok"
$t.Collapse(0)

$t.InsertAfter("""%.15f""")
$t.Style = "StringTok"
$t.Collapse(0)

$t.InsertAfter(", ")
$t.Style = "NormalTok"
$t.Collapse(0)

$t.InsertAfter("log")
$t.Style = "KeywordTok"
$t.Collapse(0)

$t.InsertAfter("(")
$t.Style = "NormalTok"
$t.Collapse(0)

$t.InsertAfter("4")
$t.Style = "DecValTok"
$t.Collapse(0)

$t.InsertAfter("))")
$t.Style = "NormalTok"
$t.Collapse(0)

# --- paragraph 2: ## [1] "1.386294361119891" ------------------------------
$o = $outPar.Range
$o.Collapse(1)
$o.InsertAfter("## [1] ""1.386294361119891""")
$o.Style = "VerbatimChar"
